$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.917.66"
$ws.Range("E2").Value = "  +3.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.415.91"
$ws.Range("E3").Value = "  +3.79%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.48"
$ws.Range("E5").Value = "  +3.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.78"
$ws.Range("E6").Value = "  +8.77%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.415.90"
$ws.Range("E8").Value = "  +3.69%  "

$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  +10.45%  "

$ws.Range("E12").Value = "  +7.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.000.82"
$ws.Range("E13").Value = "  +3.67%  "

$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  +8.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.418.00"
$ws.Range("E16").Value = "  +3.59%  "

$ws.Range("E17").Value = "  +6.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.892.59"
$ws.Range("E18").Value = "  +3.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.17"
$ws.Range("E19").Value = "  +7.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  +5.84%  "

$ws.Range("E21").Value = "  +8.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.09"
$ws.Range("E22").Value = "  +12.22%  "

$ws.Range("E23").Value = "  +4.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.553.54"
$ws.Range("E24").Value = "  +3.81%  "

$ws.Range("E25").Value = "  +20.15%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.61"
$ws.Range("E27").Value = "  +4.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.60"
$ws.Range("E28").Value = "  +11.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.68"
$ws.Range("E29").Value = "  +5.94%  "

$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.31"
$ws.Range("E31").Value = "  +7.30%  "

$ws.Range("E32").Value = "  +6.35%  "

$ws.Range("E33").Value = "  +4.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.448.00"
$ws.Range("E34").Value = "  +3.88%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.61"
$ws.Range("E36").Value = "  +4.53%  "

$ws.Range("E37").Value = "  +4.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.01"
$ws.Range("E38").Value = "  +4.44%  "

$ws.Range("E39").Value = "  +6.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.83"
$ws.Range("E40").Value = "  +4.05%  "

$ws.Range("E41").Value = "  +6.67%  "

$ws.Range("E42").Value = "  +16.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.791"
$ws.Range("E43").Value = "  +7.19%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.23"
$ws.Range("E45").Value = "  +7.28%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.48"
$ws.Range("E46").Value = "  +5.15%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.25"
$ws.Range("E47").Value = "  +13.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.69"
$ws.Range("E48").Value = "  +3.46%  "

$ws.Range("E49").Value = "  +4.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.17"
$ws.Range("E50").Value = "  +7.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.386.31"
$ws.Range("E51").Value = "  +10.54%  "
